# Rename "Form" sheet to "DIA", fix a typo in the filename list, and add
# a new "DDA" sheet (sibling of DIA) with the corresponding DDA file names.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "DIA"

# Fix typo: double underscore -> single underscore
$ws1.Range("A15").Value = '2022MH003_NIAR_005_02_30pto_DIA'

# Leave the cursor where the author left it on the DIA sheet
$ws1.Range("A16").Select()

# Add the new DDA sheet right after DIA
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "DDA"

# Mirror the sheet-scoped defined name onto the new sheet (Excel/the CSV
# query-table wizard creates one such name per imported sheet)
$ws2.Names.Add("ListPTM_quantification_of_histones___A._Pilot_experiment_0_2022_06_29_11_39_14", "=DDA!`$A`$1:`$E`$2")

$ws2.Range("A1").Value = 'Proteomic_Facility_FileName'
$ws2.Range("B1").Value = 'Short_Name'
$ws2.Range("C1").Value = 'Name'
$ws2.Range("D1").Value = 'Condition'
$ws2.Range("E1").Value = 'Replicate'
$ws2.Range("A2").Value = '2022MZ006_IVMO_001_02_25pto_DDA'
$ws2.Range("B2").Value = 'WT_1'
$ws2.Range("C2").Value = 'ZWT'
$ws2.Range("D2").Value = 'WT'
$ws2.Range("E2").Value = 1
$ws2.Range("A3").Value = '2022MZ006_IVMO_002_02_25pto_DDA'
$ws2.Range("B3").Value = 'WT_2'
$ws2.Range("C3").Value = 'ZA1'
$ws2.Range("D3").Value = 'WT'
$ws2.Range("E3").Value = 2
$ws2.Range("A4").Value = '2022MZ006_IVMO_003_02_25pto_DDA'
$ws2.Range("B4").Value = 'WT_3'
$ws2.Range("C4").Value = 'ZA2'
$ws2.Range("D4").Value = 'WT'
$ws2.Range("E4").Value = 3
$ws2.Range("A5").Value = '2022MZ006_IVMO_004_02_25pto_DDA'
$ws2.Range("B5").Value = 'CSex4_1'
$ws2.Range("C5").Value = 'ZA3'
$ws2.Range("D5").Value = 'CSexon4'
$ws2.Range("E5").Value = 1
$ws2.Range("A6").Value = '2022MZ006_IVMO_005_02_25pto_DDA'
$ws2.Range("B6").Value = 'CSex4_2'
$ws2.Range("C6").Value = 'ZA3'
$ws2.Range("D6").Value = 'CSexon4'
$ws2.Range("E6").Value = 2
$ws2.Range("A7").Value = '2022MZ006_IVMO_006_02_25pto_DDA'
$ws2.Range("B7").Value = 'CSex4_3'
$ws2.Range("C7").Value = 'ZA3'
$ws2.Range("D7").Value = 'CSexon4'
$ws2.Range("E7").Value = 3
$ws2.Range("A8").Value = '2022MZ006_IVMO_007_02_25pto_DDA'
$ws2.Range("B8").Value = 'Dex4_1'
$ws2.Range("C8").Value = 'ZO6'
$ws2.Range("D8").Value = 'Dexon4'
$ws2.Range("E8").Value = 1
$ws2.Range("A9").Value = '2022MZ006_IVMO_008_02_25pto_DDA'
$ws2.Range("B9").Value = 'Dex4_2'
$ws2.Range("C9").Value = 'ZO7'
$ws2.Range("D9").Value = 'Dexon4'
$ws2.Range("E9").Value = 2
$ws2.Range("A10").Value = '2022MZ006_IVMO_009_02_25pto_DDA'
$ws2.Range("B10").Value = 'Dex4_3'
$ws2.Range("C10").Value = 'ZO8'
$ws2.Range("D10").Value = 'Dexon4'
$ws2.Range("E10").Value = 3
$ws2.Range("A11").Value = '2022MH003_NIAR_001_01_30pto_DDA                                                                                                  '
$ws2.Range("B11").Value = 'KO_1'
$ws2.Range("C11").Value = 'ZKO1'
$ws2.Range("D11").Value = 'KO'
$ws2.Range("E11").Value = 1
$ws2.Range("A12").Value = '2022MH003_NIAR_002_01_30pto_DDA'
$ws2.Range("B12").Value = 'KO_2'
$ws2.Range("C12").Value = 'ZKO3'
$ws2.Range("D12").Value = 'KO'
$ws2.Range("E12").Value = 2
$ws2.Range("A13").Value = '2022MH003_NIAR_003_01_30pto_DDA'
$ws2.Range("B13").Value = 'KO_3'
$ws2.Range("C13").Value = 'ZKO4'
$ws2.Range("D13").Value = 'KO'
$ws2.Range("E13").Value = 3
$ws2.Range("A14").Value = '2022MH003_NIAR_004_01_30pto_DDA'
$ws2.Range("B14").Value = 'KOrL_1'
$ws2.Range("C14").Value = 'ZRL1'
$ws2.Range("D14").Value = 'KO_L'
$ws2.Range("E14").Value = 1
$ws2.Range("A15").Value = '2022MH003_NIAR_005_01_30pto_DDA'
$ws2.Range("B15").Value = 'KOrL_2'
$ws2.Range("C15").Value = 'ZRL2'
$ws2.Range("D15").Value = 'KO_L'
$ws2.Range("E15").Value = 2
$ws2.Range("A16").Value = '2022MH003_NIAR_006_01_30pto_DDA'
$ws2.Range("B16").Value = 'KOrL_3'
$ws2.Range("C16").Value = 'ZRL3'
$ws2.Range("D16").Value = 'KO_L'
$ws2.Range("E16").Value = 3
$ws2.Range("A17").Value = '2022MH003_NIAR_007_01_30pto_DDA'
$ws2.Range("B17").Value = 'KOrS_1'
$ws2.Range("C17").Value = 'ZRS4'
$ws2.Range("D17").Value = 'KO_S'
$ws2.Range("E17").Value = 1
$ws2.Range("A18").Value = '2022MH003_NIAR_008_01_30pto_DDA'
$ws2.Range("B18").Value = 'KOrS_2'
$ws2.Range("C18").Value = 'ZRS5'
$ws2.Range("D18").Value = 'KO_S'
$ws2.Range("E18").Value = 2
$ws2.Range("A19").Value = '2022MH003_NIAR_009_01_30pto_DDA'
$ws2.Range("B19").Value = 'KOrS_3'
$ws2.Range("C19").Value = 'ZRS6'
$ws2.Range("D19").Value = 'KO_S'
$ws2.Range("E19").Value = 3

$ws2.Columns.Item(1).ColumnWidth = 39.25
$ws2.Columns.Item(2).ColumnWidth = 10

$ws2.Range("A20").Select()
